$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.0008521033157308618,
    0.0076424118325185255,
    0.017274196866825384,
    0.0015168823704226227,
    0.0074830425169901404,
    0.003088629587174689,
    0.008685366435123952,
    0.009147014927810272,
    0.01993697074201439,
    0.010130271352385804,
    0.018679781977551133,
    0.002662120722256513
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

$rng = $ws.Range("C2:C13")
$rng.Style = "Percent"
$rng.NumberFormat = "0.0%"

$ws.Range("C15").Select()
